# LBCB 3 calibration fixes.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Row 1 header re-layout on "SW Cmd Cals":
#   Before:  B1:C1=Full Extension | D1:E1=Full Retraction | F1:H1=Expected(old text) | I1:K1=Mid-Point
#   After:   B1:C1=Full Extension | D1:F1=Full Retraction                            | G1:I1=Mid-Point
# ---------------------------------------------------------------
$ws2.Range("I1:K1").UnMerge()
$ws2.Range("F1:H1").UnMerge()
$ws2.Range("G1").Value = "Mid-Point"
$ws2.Range("F1").Value = ""
$ws2.Range("I1").Value = ""
$ws2.Range("G1:I1").Merge()
$ws2.Range("D1:F1").Merge()

# ---------------------------------------------------------------
# Row 2 sub-headers shift left by one (old F/G/H/I/J/K -> new F/G/H/I/J)
# and a new "Expected" label column (J2) appears.
# ---------------------------------------------------------------
$ws2.Range("F2").Value = "Inches"
$ws2.Range("G2").Value = "%FS"
$ws2.Range("H2").Value = "LVDT Volts"
$ws2.Range("I2").Value = "Inches"
$ws2.Range("J2").Value = "Expected"
$ws2.Range("K2").Value = ""

# ---------------------------------------------------------------
# New small table in columns M:O (rows 2-6)
# ---------------------------------------------------------------
$ws2.Range("M2").Value = 5.937
$ws2.Range("N2").Formula = "=M2-M3"
$ws2.Range("O2").Value = 48.5

$ws2.Range("M3").Value = 5.95
$ws2.Range("O3").Value = 48.6

$ws2.Range("O4").Value = 48.65

$ws2.Range("M5").Value = 5.97
$ws2.Range("O5").Value = -48.8

$ws2.Range("M6").Value = 6.02
$ws2.Range("O6").Value = -49.5

$ws2.Range("O7").Value = -49

# ---------------------------------------------------------------
# Updated calibration data values (rows 3-8)
# ---------------------------------------------------------------
$ws2.Range("B3").Value = -99.35
$ws2.Range("C3").Value = 10.3202
$ws2.Range("D3").Value = 0.9
$ws2.Range("E3").Value = -0.0929897

$ws2.Range("B4").Value = -98.39
$ws2.Range("C4").Value = 10.237
$ws2.Range("D4").Value = 0.45
$ws2.Range("E4").Value = -0.050712

$ws2.Range("B5").Value = -99
$ws2.Range("C5").Value = 10.12
$ws2.Range("D5").Value = 1.31
$ws2.Range("E5").Value = -0.1043
$ws2.Range("F5").Formula = "=1.03"
$ws2.Range("G5").Value = 0
$ws2.Range("H5").Value = 4.96152
$ws2.Range("I5").Value = 6.03

$ws2.Range("D6").Value = 0.09
$ws2.Range("E6").Value = -0.0163535
$ws2.Range("F6").Value = 1.03
$ws2.Range("G6").Value = -48.58
$ws2.Range("H6").Value = 4.87058
$ws2.Range("I6").Value = 6.03
$ws2.Range("J6").Value = 6.03

$ws2.Range("D7").Value = 1.15
$ws2.Range("E7").Value = -0.100279
$ws2.Range("F7").Value = 1.03
$ws2.Range("G7").Value = -48.4
$ws2.Range("H7").Value = 5.03454
$ws2.Range("I7").Value = 6.03

$ws2.Range("D8").Value = 0.26
$ws2.Range("E8").Value = -0.0217311
$ws2.Range("F8").Value = 1
$ws2.Range("G8").Value = -49
$ws2.Range("H8").Value = 5.06246
$ws2.Range("I8").Value = 6
$ws2.Range("J8").Value = 6

# ---------------------------------------------------------------
# Column widths: D narrower, H wider (to fit the new "Expected" column)
# ---------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 9.166666666666666
$ws2.Columns.Item(8).ColumnWidth = 14.333333333333334

# ---------------------------------------------------------------
# Sheet tab / selection swap: "SW Cmd Cals" becomes the active sheet/tab,
# "Displacement Cals" keeps a plain selection.
# ---------------------------------------------------------------
$ws1.Range("C7").Select()
$ws2.Activate()
$ws2.Range("H9").Select()
